# Kadastro App: Kayıt silindi: 11161048
# Delete the record row whose "Kayıt No" (column A) equals 11161048
# from the master "Kayitlar" sheet and from the filtered district
# sheet ("Merkez İlçe") that mirrors it.

$wb = $excel.ActiveWorkbook
$recordId = "11161048"

$sheetNames = @("Kayitlar", "Merkez İlçe")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $found = $ws.Columns.Item(1).Find($recordId)
    if ($found) {
        $ws.Rows.Item($found.Row).Delete()
    }
}
